# Insert a new weekly record as row 9, pushing the existing rows 9-20
# (Feria Lagunitas de Puerto Montt / Arándano (blue) price history) down
# to rows 10-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:20 down to 10:21, carrying their formatting (incl. the
# date number-format on column D) along with them.
$ws.Rows("9:9").Insert()

# Populate the newly-opened row 9 with the new weekly entry.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 44880
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101001
$ws.Range("J9").Value = "Arándano (blue)"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 7500
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 7750
$ws.Range("Q9").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 5167
$ws.Range("T9").Value = 1.5
